$wb = $excel.ActiveWorkbook

# --- Sheet1 edits ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the stray formula in G15 (0.197*360)
$ws1.Range("G15").ClearContents()

# Numeric helper column D (rows 21-32), no new shared strings created here
$ws1.Range("D20").Value = "Hour"
$ws1.Range("D21").Value = 0
$ws1.Range("D22").Value = 1
$ws1.Range("D23").Value = 2
$ws1.Range("D24").Value = 3
$ws1.Range("D25").Value = 4
$ws1.Range("D26").Value = 5
$ws1.Range("D27").Value = 6
$ws1.Range("D28").Value = 7
$ws1.Range("D29").Value = 8
$ws1.Range("D30").Value = 9
$ws1.Range("D31").Value = 10
$ws1.Range("D32").Value = 11

# Text cells set in the exact order needed to reproduce the shared-string table
$ws1.Range("E21").Value = "12 и 11"
$ws1.Range("E22").Value = "1 и 12"
$ws1.Range("E23").Value = "2 и 1"
$ws1.Range("E24").Value = "3 и 2"
$ws1.Range("E31").Value = "10 и 9"
$ws1.Range("E32").Value = "11 и 10"
$ws1.Range("G20").Value = "Hyperminute"
$ws1.Range("H21").Value = "12 и 11.5"
$ws1.Range("H22").Value = "0.5 и 0"
$ws1.Range("I22").Value = "!"
$ws1.Range("H23").Value = "1 и 0.5"
$ws1.Range("A32").Value = "F CPU"
$ws1.Range("C32").Value = "MHz"
$ws1.Range("A39").Value = "Timer input freq"
$ws1.Range("A33").Value = "Timer ovf freq"
$ws1.Range("A38").Value = "Divisor"
$ws1.Range("A40").Value = "ICR"
$ws1.Range("A41").Value = "OVF freq"

# Remaining numeric / formula cells
$ws1.Range("B32").Value = 1
$ws1.Range("B33").Formula = "=B32*1000000/256"
$ws1.Range("C33").Value = "Hz"

$ws1.Range("B38").Value = 1
$ws1.Range("B39").Formula = "=B32*1000000/B38"
$ws1.Range("C39").Value = "Hz"

$ws1.Range("B40").Value = 255
$ws1.Range("B41").Formula = "=B39/B40"
$ws1.Range("C41").Value = "Hz"

$ws1.Range("G21").Value = 0
$ws1.Range("G22").Formula = "=G21+1"
$ws1.Range("G23").Formula = "=G22+1"
$ws1.Range("G24").Formula = "=G23+1"
$ws1.Range("G25").Formula = "=G24+1"
$ws1.Range("G26").Formula = "=G25+1"
$ws1.Range("G27").Formula = "=G26+1"
$ws1.Range("G28").Formula = "=G27+1"
$ws1.Range("G29").Formula = "=G28+1"
$ws1.Range("G30").Formula = "=G29+1"
$ws1.Range("G31").Formula = "=G30+1"
$ws1.Range("G32").Formula = "=G31+1"
$ws1.Range("G33").Formula = "=G32+1"
$ws1.Range("G34").Formula = "=G33+1"
$ws1.Range("G35").Formula = "=G34+1"
$ws1.Range("G36").Formula = "=G35+1"
$ws1.Range("G37").Formula = "=G36+1"
$ws1.Range("G38").Formula = "=G37+1"
$ws1.Range("G39").Formula = "=G38+1"
$ws1.Range("G40").Formula = "=G39+1"
$ws1.Range("G41").Formula = "=G40+1"
$ws1.Range("G42").Formula = "=G41+1"
$ws1.Range("G43").Formula = "=G42+1"
$ws1.Range("G44").Formula = "=G43+1"

# Column E width
$ws1.Columns.Item(5).ColumnWidth = 10.28515625

# Selection / view state for Sheet1
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws1.Range("C44").Select()

# --- Sheet "Money" view state: no longer the tab-selected sheet ---
$wsMoney = $wb.Worksheets.Item("Money")
$wsMoney.Range("A3:A5").Select()

# Sheet1 must be the active/selected tab at the end
$ws1.Activate()
